$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2:A11 values from 4 to 1
$ws.Range("A2:A11").Value = 1

# Update selection to A2:A11 with active cell A2
$ws.Range("A2:A11").Select() | Out-Null
